$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.584.79"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "1.756.91"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "324.89"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").Value = "0.4577"
$ws.Range("E7").Value = "  +3.19%  "

$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("D9").Value = "0.07461"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "41.62"
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("D11").Value = "1.087"
$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").Value = "20.81"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").Value = "6.013"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").Value = "7.167"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "1.755.87"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").Value = "93.56"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "0.00001054"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D19").Value = "0.06430"
$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").Value = "17.06"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("D22").Value = "5.746"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").Value = "27.625.37"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").Value = "11.23"
$ws.Range("E24").Value = "  +0.47%  "

$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").Value = "165.35"
$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").Value = "20.13"
$ws.Range("E27").Value = "  -1.20%  "

$ws.Range("D28").Value = "1.955.41"
$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("D29").Value = "2.143"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("D30").Value = "125.67"
$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("D31").Value = "1.079"
$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("D32").Value = "0.09207"
$ws.Range("E32").Value = "  +2.48%  "

$ws.Range("D33").Value = "3.666"
$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").Value = "5.512"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("E35").Value = "  -1.67%  "

$ws.Range("D36").Value = "0.02277"
$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("D37").Value = "0.06011"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("D38").Value = "0.2085"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "0.6269"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").Value = "4.927"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").Value = "1.181"
$ws.Range("E41").Value = "  -2.34%  "

$ws.Range("D42").Value = "1.385"
$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("D43").Value = "7.747"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").Value = "13.25"
$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("D45").Value = "3.718"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").Value = "0.5856"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").Value = "121.92"
$ws.Range("E47").Value = "  +0.70%  "

$ws.Range("D48").Value = "1.939"
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").Value = "0.06911"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("D50").Value = "1.134"
$ws.Range("E50").Value = "  -1.60%  "

$ws.Range("D51").Value = "72.12"
$ws.Range("E51").Value = "  +0.01%  "
